$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 68004387.54000001
$ws.Range("P2").Value = 146.8657919132
$ws.Range("Q2").Value = 782790548.14
$ws.Range("R2").Value = 1690.5549467255
$ws.Range("S2").Value = 156208401.52
$ws.Range("T2").Value = 337.3557416313
$ws.Range("U2").Value = -91905660.54000001
$ws.Range("V2").Value = -198.4842170452
$ws.Range("W2").Value = 3903093.15
$ws.Range("X2").Value = 8.4293218
$ws.Range("Y2").Value = 47146429.31
$ws.Range("Z2").Value = 101.8198667317
$ws.Range("AA2").Value = -21705582.5
$ws.Range("AB2").Value = -46.8764983866
$ws.Range("AC2").Value = -46303762.54
$ws.Range("AD2").Value = -310.2546130626
